# Add localized strings for the service-worker "update available" modal.
#
# Two new rows are appended to the KeyValuePairs table (Tabelle2):
#   row 48: service_worker-update_headline        | Update Available      | Update verfügbar
#   row 49: service_worker-update_confirm_btn_txt  | Update Now & Refresh  | Update installieren
#
# Columns: A = Key, B = String EN, C = String DE

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the Key column first for both new rows, then the EN/DE strings for
# each row in turn - this mirrors how the rows were authored and keeps the
# shared-string table's insertion order aligned with the source edit.
$ws.Range("A48").Value = "service_worker-update_headline"
$ws.Range("A49").Value = "service_worker-update_confirm_btn_txt"

$ws.Range("B48").Value = "Update Available"
$ws.Range("C48").Value = "Update verfügbar"

$ws.Range("B49").Value = "Update Now & Refresh"
$ws.Range("C49").Value = "Update installieren"

# Match the formatting of the existing table body (light themed fill) by
# copying the format from an existing data row onto the two new rows.
$ws.Range("A2:C2").Copy()
$ws.Range("A48:C49").PasteSpecial(-4122)

# Grow the table ("Tabelle2") so the new rows become part of it, which also
# extends its AutoFilter range.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C49"))

# Reflect the selection left behind after entering the new data.
$ws.Range("D44").Select()
